# Applies scheduled market-data refresh to the Ixion profit-tracking workbook.
# For each "Leve" table (one per crafting job sheet) this updates the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) with freshly
# pulled price data, matching the upstream diff cell-for-cell.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2083963.5
$ws.Range("J17").Value = 2083963.5
$ws.Range("L17").Value = 6251890.5
$ws.Range("N17").Value = -6252226.5
$ws.Range("H97").Value = 4433.3335
$ws.Range("J97").Value = 4433.3335
$ws.Range("L97").Value = 13300.0005
$ws.Range("N97").Value = -14292.0005
$ws.Range("H129").Value = 970.4918
$ws.Range("I129").Value = 865.6667
$ws.Range("J129").Value = 981.92725
$ws.Range("K129").Value = 2597.0001
$ws.Range("L129").Value = 2945.78175
$ws.Range("M129").Value = 2402.9999
$ws.Range("N129").Value = -12945.78175
$ws.Range("H131").Value = 2985
$ws.Range("I131").Value = 3095
$ws.Range("J131").Value = 2875
$ws.Range("K131").Value = 9285
$ws.Range("L131").Value = 8625
$ws.Range("M131").Value = -4245
$ws.Range("N131").Value = -18705
$ws.Range("H132").Value = 1198.7637
$ws.Range("I132").Value = 932.09753
$ws.Range("J132").Value = 1979.7142
$ws.Range("K132").Value = 2796.29259
$ws.Range("L132").Value = 5939.142599999999
$ws.Range("M132").Value = -266.29259
$ws.Range("N132").Value = -10999.1426
$ws.Range("H137").Value = 1328.4445
$ws.Range("I137").Value = 1066.3334
$ws.Range("J137").Value = 2009.9333
$ws.Range("K137").Value = 3199.0002
$ws.Range("L137").Value = 6029.7999
$ws.Range("M137").Value = -649.0001999999999
$ws.Range("N137").Value = -11129.7999
$ws.Range("H138").Value = 2274.0942
$ws.Range("I138").Value = 992.1579
$ws.Range("J138").Value = 3310.5532
$ws.Range("K138").Value = 2976.4737
$ws.Range("L138").Value = 9931.659599999999
$ws.Range("M138").Value = 2163.5263
$ws.Range("N138").Value = -20211.6596
$ws.Range("H141").Value = 1638.8334
$ws.Range("I141").Value = 1174.5714
$ws.Range("J141").Value = 3263.75
$ws.Range("K141").Value = 3523.7142
$ws.Range("L141").Value = 9791.25
$ws.Range("M141").Value = 1656.2858
$ws.Range("N141").Value = -20151.25

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 925.3125
$ws.Range("I2").Value = 1022.46155
$ws.Range("J2").Value = 504.33334
$ws.Range("K2").Value = 1022.46155
$ws.Range("L2").Value = 504.33334
$ws.Range("M2").Value = -909.46155
$ws.Range("N2").Value = -730.33334
$ws.Range("H32").Value = 1361.98
$ws.Range("I32").Value = 1145.8372
$ws.Range("J32").Value = 2689.7144
$ws.Range("K32").Value = 1145.8372
$ws.Range("L32").Value = 2689.7144
$ws.Range("M32").Value = -858.8371999999999
$ws.Range("N32").Value = -3263.7144
$ws.Range("H61").Value = 2034.2858
$ws.Range("I61").Value = 2167.6667
$ws.Range("J61").Value = 1856.4445
$ws.Range("K61").Value = 2167.6667
$ws.Range("L61").Value = 1856.4445
$ws.Range("M61").Value = -1955.6667
$ws.Range("N61").Value = -2280.4445
$ws.Range("H74").Value = 1408.75
$ws.Range("I74").Value = 1219.8182
$ws.Range("J74").Value = 1824.4
$ws.Range("K74").Value = 1219.8182
$ws.Range("L74").Value = 1824.4
$ws.Range("M74").Value = -345.8181999999999
$ws.Range("N74").Value = -3572.4
$ws.Range("H77").Value = 1408.75
$ws.Range("I77").Value = 1219.8182
$ws.Range("J77").Value = 1824.4
$ws.Range("K77").Value = 6099.090999999999
$ws.Range("L77").Value = 9122
$ws.Range("M77").Value = -1731.090999999999
$ws.Range("N77").Value = -17858
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 1000
$ws.Range("M110").Value = 1045
$ws.Range("H116").Value = 925.3125
$ws.Range("I116").Value = 1022.46155
$ws.Range("J116").Value = 504.33334
$ws.Range("K116").Value = 1022.46155
$ws.Range("L116").Value = 504.33334
$ws.Range("M116").Value = 1271.53845
$ws.Range("N116").Value = -5092.33334
$ws.Range("H122").Value = 1973272.4
$ws.Range("I122").Value = 2849694
$ws.Range("J122").Value = 1323.75
$ws.Range("K122").Value = 8549082
$ws.Range("L122").Value = 3971.25
$ws.Range("M122").Value = -8546632
$ws.Range("N122").Value = -8871.25
$ws.Range("H132").Value = 2130391
$ws.Range("I132").Value = 1896.0938
$ws.Range("J132").Value = 6671180.5
$ws.Range("K132").Value = 5688.2814
$ws.Range("L132").Value = 20013541.5
$ws.Range("M132").Value = -3158.2814
$ws.Range("N132").Value = -20018601.5
$ws.Range("H136").Value = 2034.2858
$ws.Range("I136").Value = 2167.6667
$ws.Range("J136").Value = 1856.4445
$ws.Range("K136").Value = 6503.000100000001
$ws.Range("L136").Value = 5569.333500000001
$ws.Range("M136").Value = -3953.000100000001
$ws.Range("N136").Value = -10669.3335

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 925.3125
$ws.Range("I3").Value = 1022.46155
$ws.Range("J3").Value = 504.33334
$ws.Range("K3").Value = 1022.46155
$ws.Range("L3").Value = 504.33334
$ws.Range("M3").Value = -908.46155
$ws.Range("N3").Value = -732.33334
$ws.Range("H134").Value = 1663.8889
$ws.Range("I134").Value = 1329.6875
$ws.Range("J134").Value = 2150
$ws.Range("K134").Value = 3989.0625
$ws.Range("L134").Value = 6450
$ws.Range("M134").Value = -1454.0625
$ws.Range("N134").Value = -11520

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21736.309
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 21736.309
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 21736.309
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -22326.309
$ws.Range("H34").Value = 21736.309
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 21736.309
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 21736.309
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -22140.309
$ws.Range("H58").Value = 1299.4916
$ws.Range("I58").Value = 721.6111
$ws.Range("J58").Value = 2204
$ws.Range("K58").Value = 721.6111
$ws.Range("L58").Value = 2204
$ws.Range("M58").Value = -518.6111
$ws.Range("N58").Value = -2610
$ws.Range("H99").Value = 6258779
$ws.Range("I99").Value = 8511.933999999999
$ws.Range("J99").Value = 25009580
$ws.Range("K99").Value = 8511.933999999999
$ws.Range("L99").Value = 25009580
$ws.Range("M99").Value = -7013.933999999999
$ws.Range("N99").Value = -25012576
$ws.Range("H126").Value = 6258779
$ws.Range("I126").Value = 8511.933999999999
$ws.Range("J126").Value = 25009580
$ws.Range("K126").Value = 25535.802
$ws.Range("L126").Value = 75028740
$ws.Range("M126").Value = -23065.802
$ws.Range("N126").Value = -75033680
$ws.Range("H132").Value = 2514.976
$ws.Range("I132").Value = 2056.074
$ws.Range("J132").Value = 3341
$ws.Range("K132").Value = 6168.222
$ws.Range("L132").Value = 10023
$ws.Range("M132").Value = -3638.222
$ws.Range("N132").Value = -15083
$ws.Range("H134").Value = 1793.5483
$ws.Range("I134").Value = 2145.3333
$ws.Range("J134").Value = 1306.4615
$ws.Range("K134").Value = 6435.999899999999
$ws.Range("L134").Value = 3919.3845
$ws.Range("M134").Value = -3900.999899999999
$ws.Range("N134").Value = -8989.3845
$ws.Range("H136").Value = 1299.4916
$ws.Range("I136").Value = 721.6111
$ws.Range("J136").Value = 2204
$ws.Range("K136").Value = 2164.8333
$ws.Range("L136").Value = 6612
$ws.Range("M136").Value = 385.1667000000002
$ws.Range("N136").Value = -11712

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 193099.47
$ws.Range("I5").Value = 270.2414
$ws.Range("J5").Value = 436231.97
$ws.Range("K5").Value = 810.7242
$ws.Range("L5").Value = 1308695.91
$ws.Range("M5").Value = -698.7242
$ws.Range("N5").Value = -1308919.91
$ws.Range("H113").Value = 238656.97
$ws.Range("I113").Value = 581.8182
$ws.Range("J113").Value = 500539.66
$ws.Range("K113").Value = 1745.4546
$ws.Range("L113").Value = 1501618.98
$ws.Range("M113").Value = 424.5454
$ws.Range("N113").Value = -1505958.98
$ws.Range("H135").Value = 193099.47
$ws.Range("I135").Value = 270.2414
$ws.Range("J135").Value = 436231.97
$ws.Range("K135").Value = 2432.1726
$ws.Range("L135").Value = 3926087.73
$ws.Range("M135").Value = 102.8274000000001
$ws.Range("N135").Value = -3931157.73

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 14333333
$ws.Range("I24").Value = 14333333
$ws.Range("K24").Value = 14333333
$ws.Range("M24").Value = -14333160
$ws.Range("H80").Value = 3148.5
$ws.Range("I80").Value = 2347
$ws.Range("K80").Value = 2347
$ws.Range("M80").Value = -1349
$ws.Range("H83").Value = 3148.5
$ws.Range("I83").Value = 2347
$ws.Range("K83").Value = 11735
$ws.Range("M83").Value = -6743
$ws.Range("H122").Value = 62638420
$ws.Range("I122").Value = 106483464
$ws.Range("J122").Value = 2646.1428
$ws.Range("K122").Value = 319450392
$ws.Range("L122").Value = 7938.428400000001
$ws.Range("M122").Value = -319447942
$ws.Range("N122").Value = -12838.4284
$ws.Range("H123").Value = 22198.6
$ws.Range("J123").Value = 22198.6
$ws.Range("L123").Value = 22198.6
$ws.Range("N123").Value = -27098.6
$ws.Range("H132").Value = 2469.3125
$ws.Range("I132").Value = 1998.7059
$ws.Range("J132").Value = 3002.6667
$ws.Range("K132").Value = 5996.1177
$ws.Range("L132").Value = 9008.000100000001
$ws.Range("M132").Value = -3466.1177
$ws.Range("N132").Value = -14068.0001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 127238
$ws.Range("I7").Value = 144772
$ws.Range("K7").Value = 144772
$ws.Range("M7").Value = -144660
$ws.Range("H40").Value = 71432104
$ws.Range("I40").Value = 90911770
$ws.Range("J40").Value = 6665
$ws.Range("K40").Value = 90911770
$ws.Range("L40").Value = 6665
$ws.Range("M40").Value = -90911634
$ws.Range("N40").Value = -6937
$ws.Range("H122").Value = 2548983.8
$ws.Range("I122").Value = 3110568.5
$ws.Range("J122").Value = 1113822.2
$ws.Range("K122").Value = 9331705.5
$ws.Range("L122").Value = 3341466.6
$ws.Range("M122").Value = -9329255.5
$ws.Range("N122").Value = -3346366.6
$ws.Range("H126").Value = 127238
$ws.Range("I126").Value = 144772
$ws.Range("K126").Value = 434316
$ws.Range("M126").Value = -431846
$ws.Range("H136").Value = 5797.8096
$ws.Range("I136").Value = 3652.3137
$ws.Range("J136").Value = 14916.167
$ws.Range("K136").Value = 10956.9411
$ws.Range("L136").Value = 44748.501
$ws.Range("M136").Value = -8406.9411
$ws.Range("N136").Value = -49848.501

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2299.923
$ws.Range("I122").Value = 2405.158
$ws.Range("J122").Value = 2014.2858
$ws.Range("K122").Value = 7215.474
$ws.Range("L122").Value = 6042.857400000001
$ws.Range("M122").Value = -4765.474
$ws.Range("N122").Value = -10942.8574
$ws.Range("H136").Value = 7578194.5
$ws.Range("I136").Value = 2564.8572
$ws.Range("K136").Value = 7694.571599999999
$ws.Range("M136").Value = -5144.571599999999
